$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.382.93"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "2.647.04"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.83"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.66"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").Value = "2.645.85"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("E10").Value = "  +7.97%  "
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.27"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.08"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("E15").Value = "  +2.92%  "
$ws.Range("D16").Value = "3.128.54"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").Value = "68.299.04"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "2.646.38"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "363.86"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.49"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("E22").Value = "  +3.35%  "
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.83"
$ws.Range("E25").Value = "  +3.88%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.82"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  +2.21%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "573.73"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.09"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.87"
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("E35").Value = "  +2.70%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.60"
$ws.Range("E37").Value = "  +5.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.19"
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").Value = "0.0₆0338"
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("E44").Value = "  +1.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.74"
$ws.Range("E45").Value = "  +2.08%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.66"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.86"
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("E49").Value = "  +1.92%  "
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.94"
$ws.Range("E51").Value = "  +0.29%  "
